# Optimize temperature chunk ticks
# Shrink the heat_capacity values (column B) on the block_temperature sheet
# so each chunk needs far fewer ticks to transition, and move the active
# selection to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("block_temperature")

$newHeatCapacity = @{
    2  = 5
    3  = 5
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 3
    9  = 5
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 2
    15 = 2
    16 = 2
    17 = 2
    18 = 2
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 10
}

foreach ($row in $newHeatCapacity.Keys) {
    $ws.Cells.Item($row, 2).Value = $newHeatCapacity[$row]
}

$ws.Range("C4").Select()
